# Daily attendance processing - 2026-01-05 23:34:21
# For every row in column G ("Recorded By"), rotate the comma-separated
# list of recorders one position to the left (the first name moves to
# the end of the list). Lists with a single entry are left unchanged.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$usedRange = $ws.UsedRange
$lastRow = $usedRange.Row + $usedRange.Rows.Count - 1

for ($r = 2; $r -le $lastRow; $r++) {
    $cell = $ws.Cells.Item($r, 7)
    $value = $cell.Value2

    if ($value -ne $null -and $value -ne "") {
        $parts = $value -split ", "
        if ($parts.Count -gt 1) {
            $first = $parts[0]
            $rest = $parts[1..($parts.Count - 1)]
            $newParts = $rest + $first
            $cell.Value = $newParts -join ", "
        }
    }
}
